# Weekly update: insert a new data row at row 13 (shifting the existing
# rows 13-37 down to 14-38) and populate it with the new week's record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = 45028
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 100112043
$ws.Range("G13").Value = "Pepino dulce"
$ws.Range("H13").Value = "Cultivar IV Región"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 15000
$ws.Range("N13").Value = "`$/bandeja 18 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 833
$ws.Range("Q13").Value = 18
$ws.Range("R13").Value = "Hortaliza"
